$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "23.447.99"
$ws.Range("E2").Value = "  -1.27%  "
Set-TextValue $ws.Range("D3") "1.645.74"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  +0.13%  "
Set-TextValue $ws.Range("D6") "298.37"
$ws.Range("E6").Value = "  -1.87%  "
Set-TextValue $ws.Range("D7") "0.3781"
$ws.Range("E7").Value = "  -1.26%  "
Set-TextValue $ws.Range("D8") "0.3551"
$ws.Range("E8").Value = "  -1.73%  "
Set-TextValue $ws.Range("D9") "49.88"
$ws.Range("E9").Value = "  -2.94%  "
Set-TextValue $ws.Range("D10") "0.08088"
$ws.Range("E10").Value = "  -1.88%  "
Set-TextValue $ws.Range("D11") "1.216"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("E12").Value = "  +0.17%  "
Set-TextValue $ws.Range("D13") "22.00"
$ws.Range("E13").Value = "  -3.39%  "
Set-TextValue $ws.Range("D14") "6.383"
$ws.Range("E14").Value = "  -2.55%  "
Set-TextValue $ws.Range("D15") "7.333"
$ws.Range("E15").Value = "  -1.16%  "
Set-TextValue $ws.Range("D16") "0.00001194"
$ws.Range("E16").Value = "  -3.41%  "
Set-TextValue $ws.Range("D17") "1.653.60"
$ws.Range("E17").Value = "  +0.05%  "
Set-TextValue $ws.Range("D18") "97.36"
$ws.Range("E18").Value = "  -0.28%  "
Set-TextValue $ws.Range("D19") "0.06952"
$ws.Range("E19").Value = "  -0.46%  "
Set-TextValue $ws.Range("D20") "6.775"
$ws.Range("E20").Value = "  -0.19%  "
Set-TextValue $ws.Range("D21") "17.28"
$ws.Range("E21").Value = "  -2.50%  "
Set-TextValue $ws.Range("D22") "1.000"
$ws.Range("E22").Value = "  +0.09%  "
Set-TextValue $ws.Range("D23") "12.40"
$ws.Range("E23").Value = "  -1.96%  "
Set-TextValue $ws.Range("D24") "23.454.44"
$ws.Range("E24").Value = "  -1.20%  "
Set-TextValue $ws.Range("D25") "2.486"
$ws.Range("E25").Value = "  -2.00%  "
Set-TextValue $ws.Range("D26") "2.891"
$ws.Range("E26").Value = "  -6.10%  "
Set-TextValue $ws.Range("D27") "20.87"
$ws.Range("E27").Value = "  -2.24%  "
Set-TextValue $ws.Range("D28") "153.16"
$ws.Range("E28").Value = "  +0.98%  "
Set-TextValue $ws.Range("D29") "5.201"
$ws.Range("E29").Value = "  -1.60%  "
Set-TextValue $ws.Range("D30") "132.62"
$ws.Range("E30").Value = "  -1.85%  "
Set-TextValue $ws.Range("D31") "1.829.84"
$ws.Range("E31").Value = "  -0.29%  "
Set-TextValue $ws.Range("D32") "6.879"
$ws.Range("E32").Value = "  -0.06%  "
Set-TextValue $ws.Range("D33") "2.122"
$ws.Range("E33").Value = "  +0.67%  "
Set-TextValue $ws.Range("D34") "11.43"
$ws.Range("E34").Value = "  -3.96%  "
Set-TextValue $ws.Range("D35") "0.9970"
$ws.Range("E35").Value = "  -8.29%  "
Set-TextValue $ws.Range("D36") "0.02706"
$ws.Range("E36").Value = "  -4.99%  "
Set-TextValue $ws.Range("D37") "0.08736"
$ws.Range("E37").Value = "  -1.17%  "
Set-TextValue $ws.Range("D38") "0.2419"
$ws.Range("E38").Value = "  -4.06%  "
Set-TextValue $ws.Range("D39") "5.899"
$ws.Range("E39").Value = "  -3.73%  "
$ws.Range("E40").Value = "  +0.60%  "
Set-TextValue $ws.Range("D41") "0.06761"
$ws.Range("E41").Value = "  -4.20%  "
Set-TextValue $ws.Range("D42") "0.6858"
$ws.Range("E42").Value = "  -3.16%  "
Set-TextValue $ws.Range("D43") "1.299"
$ws.Range("E43").Value = "  -3.10%  "
Set-TextValue $ws.Range("D44") "15.41"
$ws.Range("E44").Value = "  -3.52%  "
$ws.Range("E45").Value = "  +0.16%  "
Set-TextValue $ws.Range("D46") "0.6349"
$ws.Range("E46").Value = "  -3.27%  "
Set-TextValue $ws.Range("D47") "2.244"
$ws.Range("E47").Value = "  -4.11%  "
Set-TextValue $ws.Range("D48") "3.906"
$ws.Range("E48").Value = "  -1.70%  "
Set-TextValue $ws.Range("D49") "0.07721"
$ws.Range("E49").Value = "  -3.39%  "
Set-TextValue $ws.Range("D50") "126.96"
$ws.Range("E50").Value = "  -1.44%  "
Set-TextValue $ws.Range("D51") "1.145"
$ws.Range("E51").Value = "  -4.44%  "
